# Lesson_Names.xlsx - "Finished adding all classique classes"
# Adds three new lesson-name rows into the sorted (Original Name / Abbreviation)
# table on Tabelle1, each inserted in alphabetical order by column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# 1) "Communication média / pratique des médias" / "ComMédia"
#    -> belongs right after "Communication média" (row 17), before "Communication visuelle"
$ws.Rows.Item(18).Insert()
$ws.Range("B18").Value2 = "Communication média / pratique des médias"
$ws.Range("C18").Value2 = "ComMédia"

# 2) "Littérature comparée" / "Littérature"
#    -> belongs right after "Latin" (old row 58 -> now row 59), before "Luxembourgeois"
$ws.Rows.Item(60).Insert()
$ws.Range("B60").Value2 = "Littérature comparée"
$ws.Range("C60").Value2 = "Littérature"
$ws.Range("C60").Style = "Normal"

# 3) "Pratique des médias" / "Pratique"
#    -> belongs right after "Physique / Chimie" (old row 66 -> now row 68), before "Pratique instrumentale"
$ws.Rows.Item(69).Insert()
$ws.Range("B69").Value2 = "Pratique des médias"
$ws.Range("C69").Value2 = "Pratique"

# Reflect the author's final cursor/scroll position when the file was saved.
$ws.Range("A43").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 43
$ws.Range("B90").Select() | Out-Null
